$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data rows (row 2 through 16), wiping both content and formatting
$ws.Range("A2:N16").Clear()

# Populate rows 2..17 with the refreshed dataset (reordered, new row for id 43, columns filled for ids 33 & 37).
# Cells are written column-by-column (matching the source system export order) so that the
# shared-string table is rebuilt in the same sequence as the target workbook.
# Column A
$ws.Range("A2").Value() = 30
$ws.Range("A3").Value() = 41
$ws.Range("A4").Value() = 43
$ws.Range("A5").Value() = 31
$ws.Range("A6").Value() = 32
$ws.Range("A7").Value() = 34
$ws.Range("A8").Value() = 38
$ws.Range("A9").Value() = 39
$ws.Range("A10").Value() = 40
$ws.Range("A11").Value() = 35
$ws.Range("A12").Value() = 36
$ws.Range("A13").Value() = 27
$ws.Range("A14").Value() = 33
$ws.Range("A15").Value() = 37
$ws.Range("A16").Value() = 28
$ws.Range("A17").Value() = 29

# Column B
$ws.Range("B2").Value() = "Giuseppe"
$ws.Range("B3").Value() = "Giuseppe"
$ws.Range("B4").Value() = "Giuseppe"
$ws.Range("B5").Value() = "Giuseppe"
$ws.Range("B6").Value() = "Giuseppe"
$ws.Range("B7").Value() = "giuseppe"
$ws.Range("B8").Value() = "sdsadsa"
$ws.Range("B9").Value() = "ds"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value() = "7"
$ws.Range("B10").Style = "Normal"
$ws.Range("B11").Value() = "giuseppe"
$ws.Range("B12").Value() = "sucuni"
$ws.Range("B13").Value() = "Giuseppe"
$ws.Range("B14").Value() = "Ci Piace la figa"
$ws.Range("B15").Value() = "A"
$ws.Range("B16").Value() = "Giuseppe"
$ws.Range("B17").Value() = "Giuseppe"

# Column C
$ws.Range("C2").Value() = "Cangemi"
$ws.Range("C3").Value() = "Cangemi"
$ws.Range("C4").Value() = "Cangemi"
$ws.Range("C5").Value() = "Cangemi"
$ws.Range("C6").Value() = "Cangemi"
$ws.Range("C7").Value() = "cangemi"
$ws.Range("C8").Value() = "dasdsadsa"
$ws.Range("C9").Value() = "sa"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value() = "2"
$ws.Range("C10").Style = "Normal"
$ws.Range("C11").Value() = "cangemi"
$ws.Range("C12").Value() = "npizzu"
$ws.Range("C13").Value() = "Cangemi"
$ws.Range("C14").Value() = "Tantissimo"
$ws.Range("C15").Value() = "B"
$ws.Range("C16").Value() = "Cangemi"
$ws.Range("C17").Value() = "Cangemi"

# Column D
$ws.Range("D2").Value() = "giuseppecangemi94@gmail.com"
$ws.Range("D3").Value() = "giuseppecangemi94@gmail.com"
$ws.Range("D4").Value() = "giuseppecangemi94@gmail.com"
$ws.Range("D5").Value() = "giuseppecangemi94@gmail.com"
$ws.Range("D6").Value() = "giuseppecangemi94@gmail.com"
$ws.Range("D7").Value() = "giuseppe.cangemi@prova.it"
$ws.Range("D8").Value() = "dsadadas@sdmsa.it"
$ws.Range("D9").Value() = "giuseppe.cangemi@prova.it"
$ws.Range("D10").Value() = "76h@kjjm.it"
$ws.Range("D11").Value() = "giuseppe.cangemi@prova.it"
$ws.Range("D12").Value() = "cututtiicugghiunq@sucuni.bisdum"
$ws.Range("D13").Value() = "giuseppecangemi94@gmail.com"
$ws.Range("D14").Value() = "salvatoremariazuccarello@pornhub.com"
$ws.Range("D15").Value() = "b@b.it"
$ws.Range("D16").Value() = "giuseppecangemi94@gmail.com"
$ws.Range("D17").Value() = "giuseppecangemi94@gmail.com"

# Column E
$ws.Range("E2").Value() = "derryrockfoto.jpg"
$ws.Range("E3").Value() = "derryrockfoto.jpg"
$ws.Range("E4").Value() = "derryrockpubfidelity.png"
$ws.Range("E5").Value() = "derryrockpubfidelity.png"
$ws.Range("E6").Value() = "derryrockfoto.jpg"
$ws.Range("E7").Value() = "image (31).png"
$ws.Range("E8").Value() = "image (30).png"
$ws.Range("E9").Value() = "Immagine 2024-10-31 122037.png"
$ws.Range("E10").Value() = "Immagine 2024-10-31 122037.png"
$ws.Range("E11").Value() = "image (30).png"
$ws.Range("E12").Value() = "image.jpg"
$ws.Range("E13").Value() = "derryrockfoto.jpg"
$ws.Range("E14").Value() = "1000051796.jpg"
$ws.Range("E15").Value() = "image.jpg"
$ws.Range("E16").Value() = "derryrockfoto.jpg"
$ws.Range("E17").Value() = "derryrockfoto.jpg"

# Column F
$ws.Range("F2").Value() = "SI"
$ws.Range("F3").Value() = "SI"
$ws.Range("F4").Value() = "SI"
$ws.Range("F5").Value() = "SI"
$ws.Range("F6").Value() = "SI"
$ws.Range("F8").Value() = "SI"
$ws.Range("F9").Value() = "SI"
$ws.Range("F10").Value() = "SI"
$ws.Range("F11").Value() = "SI"
$ws.Range("F12").Value() = "SI"
$ws.Range("F13").Value() = "SI"
$ws.Range("F14").Value() = "SI"
$ws.Range("F15").Value() = "SI"
$ws.Range("F16").Value() = "SI"
$ws.Range("F17").Value() = "SI"

# Column G
$ws.Range("G2").Value() = 3501
$ws.Range("G3").Value() = 7002
$ws.Range("G4").Value() = 7003
$ws.Range("G5").Value() = 7000
$ws.Range("G6").Value() = 7001
$ws.Range("G8").Value() = 7001
$ws.Range("G9").Value() = 7000
$ws.Range("G10").Value() = 7000
$ws.Range("G11").Value() = 7001
$ws.Range("G12").Value() = 666
$ws.Range("G13").Value() = 1201
$ws.Range("G14").Value() = 9999
$ws.Range("G15").Value() = 10000
$ws.Range("G16").Value() = 6201
$ws.Range("G17").Value() = 3500

# Column H
$ws.Range("H2").Value() = "SI"
$ws.Range("H3").Value() = "SI"
$ws.Range("H4").Value() = "SI"
$ws.Range("H5").Value() = "SI"
$ws.Range("H6").Value() = "SI"
$ws.Range("H8").Value() = "SI"
$ws.Range("H9").Value() = "SI"
$ws.Range("H10").Value() = "SI"
$ws.Range("H11").Value() = "SI"
$ws.Range("H12").Value() = "SI"
$ws.Range("H13").Value() = "SI"
$ws.Range("H14").Value() = "SI"
$ws.Range("H15").Value() = "SI"
$ws.Range("H16").Value() = "SI"
$ws.Range("H17").Value() = "SI"

# Column I
$ws.Range("I5").Value() = "Sì"
$ws.Range("I12").Value() = "Sì"
$ws.Range("I14").Value() = "Sì"
$ws.Range("I16").Value() = "Sì"
$ws.Range("I17").Value() = "Sì"

# Column J
$ws.Range("J3").Value() = "MOTTA SANT'ANASTASIA"
$ws.Range("J4").Value() = "MOTTA SANT'ANASTASIA"
$ws.Range("J7").Value() = "Milano"
$ws.Range("J8").Value() = "dadadas"
$ws.Range("J9").Value() = "Milano"
$ws.Range("J10").Value() = "sas"
$ws.Range("J11").Value() = "Milano"
$ws.Range("J12").Value() = "carropepe"
$ws.Range("J15").Value() = "N"

# Column K
$ws.Range("K3").NumberFormat = "YYYY-MM-DD"
$ws.Range("K3").Value() = 34681
$ws.Range("K4").NumberFormat = "YYYY-MM-DD"
$ws.Range("K4").Value() = 34315
$ws.Range("K7").NumberFormat = "YYYY-MM-DD"
$ws.Range("K7").Value() = 34680
$ws.Range("K8").NumberFormat = "YYYY-MM-DD"
$ws.Range("K8").Value() = -644653
$ws.Range("K9").NumberFormat = "YYYY-MM-DD"
$ws.Range("K9").Value() = -280525
$ws.Range("K10").NumberFormat = "YYYY-MM-DD"
$ws.Range("K10").Value() = -214781
$ws.Range("K11").NumberFormat = "YYYY-MM-DD"
$ws.Range("K11").Value() = 34680
$ws.Range("K12").NumberFormat = "YYYY-MM-DD"
$ws.Range("K12").Value() = -5662
$ws.Range("K15").NumberFormat = "YYYY-MM-DD"
$ws.Range("K15").Value() = 45573

# Column L
$ws.Range("L3").Value() = "e"
$ws.Range("L4").Value() = "s"
$ws.Range("L7").Value() = "Motta Sant'Anastasia"
$ws.Range("L8").Value() = "Motta Sant'Anastasia"
$ws.Range("L9").Value() = "SAA"
$ws.Range("L10").Value() = "ds"
$ws.Range("L11").Value() = "Motta Sant'Anastasia"
$ws.Range("L12").Value() = "misterbianco "
$ws.Range("L15").Value() = "H"

# Column M
$ws.Range("M3").Value() = "Maschile"
$ws.Range("M4").Value() = "Maschile"
$ws.Range("M7").Value() = "Maschile"
$ws.Range("M8").Value() = "Femminile"
$ws.Range("M9").Value() = "Femminile"
$ws.Range("M10").Value() = "Maschile"
$ws.Range("M11").Value() = "Maschile"
$ws.Range("M12").Value() = "Altro"
$ws.Range("M15").Value() = "Altro"

# Column N
$ws.Range("N3").Value() = "eee"
$ws.Range("N4").Value() = "s"
$ws.Range("N7").Value() = "id1234"
$ws.Range("N8").Value() = "id1234"
$ws.Range("N9").Value() = "SAD"
$ws.Range("N10").Value() = "AHHAHAHAHAHAHA"
$ws.Range("N11").Value() = "id1234"
$ws.Range("N12").NumberFormat = "@"
$ws.Range("N12").Value() = "161819273628191"
$ws.Range("N12").Style = "Normal"
$ws.Range("N15").Value() = "Jaj"

Write-Output "done"